# Updates the two-digit multiplication problems in the document.
# Each old problem string is unique within the document, so a simple
# Find/Replace (ReplaceAll) per pair is sufficient and safe.

$d = $word.ActiveDocument

$pairs = @(
    @("32×34=", "73×47="),
    @("70×98=", "54×11="),
    @("24×19=", "72×67="),
    @("14×87=", "69×56="),
    @("78×54=", "40×31="),
    @("86×86=", "51×88="),
    @("11×44=", "32×41="),
    @("70×71=", "45×41="),
    @("67×79=", "24×21="),
    @("30×69=", "49×17="),
    @("28×89=", "97×95="),
    @("29×48=", "56×39="),
    @("99×84=", "34×50="),
    @("81×71=", "83×15="),
    @("47×18=", "47×79="),
    @("58×84=", "34×29="),
    @("17×78=", "76×78="),
    @("21×54=", "23×76="),
    @("91×74=", "30×46="),
    @("42×34=", "99×92="),
    @("66×98=", "56×85="),
    @("48×64=", "74×20="),
    @("15×62=", "80×32="),
    @("65×70=", "38×74="),
    @("64×64=", "13×27=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
